$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.783.45'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -7.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.690.11'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -7.21%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.06'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -5.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.51'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.70%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.677.56'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -7.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.632'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -7.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.712'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -5.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.165'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -12.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '52.26'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -7.63%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000300'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -12.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.60'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -4.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.282.19'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -7.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.721.57'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -6.78%  '
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.127'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.16%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.35'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -5.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.98'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -7.91%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -8.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '67.825.01'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -7.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '408.16'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -7.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.63'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -4.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.50'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -7.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.07'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -8.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.81'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -9.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.72'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.72%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.81'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -6.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.93'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.55'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -8.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.04'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.72%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.85'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -8.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.72'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -7.04%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.118'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -9.04%  '
$ws.Range("B35").Value = 'InjectiveProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '44.53'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -6.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '66.10'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -7.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0₃0925'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -8.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '598.43'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -5.92%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.403'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -6.79%  '
$ws.Range("B40").Value = 'Dai'
$ws.Range("C40").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.32'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +14.31%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.136'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -7.36%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.06'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -11.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0441'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -8.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.54'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -12.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.58'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.26%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -9.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.745.17'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.52%  '
$ws.Range("B50").Value = 'WEMIXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.65'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -18.55%  '
$ws.Range("B51").Value = 'ApeXProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.11'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -8.83%  '
